$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "29.308.64", "0.7100") that must
# stay as literal text, matching the source inlineStr cells. Pre-format the data
# range as Text so Excel does not coerce the assigned strings into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.323.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.874.83"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7102"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "241.95"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.07804"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "0.3109"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "25.16"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "0.08427"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.875.89"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "5.242"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "0.7133"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "91.10"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "29.328.56"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "0.000008318"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "6.085"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "240.63"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "13.23"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "2.112.39"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "7.750"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "0.1596"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "162.39"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "9.027"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D30").Value = "4.403"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").Value = "4.326"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "0.05363"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("D34").Value = "1.949"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "0.7505"
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "0.01882"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "1.226.26"
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "6.489"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").Value = "0.8935"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "109.22"
$ws.Range("E43").Value = "  +5.16%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "72.43"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "2.009.19"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.798"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.5199"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "9.450"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "0.4326"
$ws.Range("E51").Value = "  +0.49%  "

# Strip the temporary text format now that the literal strings are committed,
# so the cells end up with no explicit style override (matching the source,
# which carries no "s" attribute on these data cells).
$ws.Range("D2:D51").ClearFormats()
